$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Insert a new column before column C ("base" and everything right of it,
#    C:Z, shifts right to D:AA) to make room for the new "aws.ses" command
#    list.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).Insert()

# ---------------------------------------------------------------------------
# 2) Populate the new column C with the aws.ses command list (header + the
#    two new commands).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 3).Value = "aws.ses"
$ws.Cells.Item(2, 3).Value = "sendMail(profile,to,subject,body)"
$ws.Cells.Item(3, 3).Value = "sendTextMail(profile,to,subject,body)"

# ---------------------------------------------------------------------------
# 3) Rewrite column A (the "target" list of command-group names) so that
#    "aws.ses" is inserted alphabetically right after "aws.s3", pushing
#    everything below it down by one row.
# ---------------------------------------------------------------------------
$targetList = @("aws.s3","aws.ses","base","csv","desktop","excel","external","image","io","jms","json","mail","number","pdf","rdbms","redis","sms","sound","ssh","step","web","webalert","webcookie","ws","ws.async","xml")
for ($i = 0; $i -lt $targetList.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $targetList[$i]
}

# ---------------------------------------------------------------------------
# 4) Update the workbook-level defined names so they point at the correct
#    (now shifted) ranges, and add the new "aws.ses" name.
# ---------------------------------------------------------------------------
$wb.Names.Item("base").RefersTo = "='#system'!`$D`$2:`$D`$36"
$wb.Names.Item("csv").RefersTo = "='#system'!`$E`$2:`$E`$5"
$wb.Names.Item("desktop").RefersTo = "='#system'!`$F`$2:`$F`$92"
$wb.Names.Item("excel").RefersTo = "='#system'!`$G`$2:`$G`$14"
$wb.Names.Item("external").RefersTo = "='#system'!`$H`$2:`$H`$3"
$wb.Names.Item("image").RefersTo = "='#system'!`$I`$2:`$I`$5"
$wb.Names.Item("io").RefersTo = "='#system'!`$J`$2:`$J`$24"
$wb.Names.Item("jms").RefersTo = "='#system'!`$K`$2:`$K`$4"
$wb.Names.Item("json").RefersTo = "='#system'!`$L`$2:`$L`$14"
$wb.Names.Item("mail").RefersTo = "='#system'!`$M`$2:`$M`$2"
$wb.Names.Item("number").RefersTo = "='#system'!`$N`$2:`$N`$15"
$wb.Names.Item("pdf").RefersTo = "='#system'!`$O`$2:`$O`$16"
$wb.Names.Item("rdbms").RefersTo = "='#system'!`$P`$2:`$P`$7"
$wb.Names.Item("redis").RefersTo = "='#system'!`$Q`$2:`$Q`$10"
$wb.Names.Item("sms").RefersTo = "='#system'!`$R`$2:`$R`$2"
$wb.Names.Item("sound").RefersTo = "='#system'!`$S`$2:`$S`$5"
$wb.Names.Item("ssh").RefersTo = "='#system'!`$T`$2:`$T`$9"
$wb.Names.Item("step").RefersTo = "='#system'!`$U`$2:`$U`$4"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$27"
$wb.Names.Item("web").RefersTo = "='#system'!`$V`$2:`$V`$117"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$W`$2:`$W`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$X`$2:`$X`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$Y`$2:`$Y`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AA`$2:`$AA`$11"

$wb.Names.Add("aws.ses", "='#system'!`$C`$2:`$C`$3")

# ---------------------------------------------------------------------------
# 5) The sheet's recorded dimension has historically run one column past the
#    last populated column (e.g. before this edit it read "A1:AA117" while
#    data only extends to column Z). Nudge that same stale margin forward by
#    one column (to "AB") to keep it consistent post-edit, without adding any
#    visible content (copy the default/no-op style of A1 onto the anchor
#    cell so no new style entries are introduced).
# ---------------------------------------------------------------------------
$ws.Cells.Item(117, 28).Style = $ws.Cells.Item(1, 1).Style
